# Campus_2_B_Park2 - rename spot ids from "A-N" to "B-N" (new DB register prefix)
# and refresh sheet view/page setup metadata to match the re-saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename spot identifiers in column A (rows 6-15): "A-1".."A-10" -> "B-1".."B-10"
for ($row = 6; $row -le 15; $row++) {
    $oldValue = $ws.Range("A$row").Text
    $newValue = $oldValue -replace '^A-', 'B-'
    $ws.Range("A$row").Value = $newValue
}

# Reset the view: scroll back to the top-left (A1) and move the active selection to A15
$win = $wb.Windows.Item(1)
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A15").Select()

# Set print page setup (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Output "Spot ids updated to B-series, view/page setup refreshed."
